$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.651.75"
$ws.Range("E2").Value = "  +4.62%  "

$ws.Range("D3").Value = "2.287.44"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.82"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.04"
$ws.Range("E7").Value = "  +2.18%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.417"
$ws.Range("E9").Value = "  +3.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0923"
$ws.Range("E10").Value = "  +4.02%  "

$ws.Range("D11").Value = "2.642.50"
$ws.Range("E11").Value = "  +3.63%  "

$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.90"
$ws.Range("E13").Value = "  +1.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.87"
$ws.Range("E14").Value = "  +9.76%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.339.99"
$ws.Range("E15").Value = "  +5.00%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.75"
$ws.Range("E16").Value = "  +3.31%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.816"
$ws.Range("E17").Value = "  +2.10%  "

$ws.Range("D18").Value = "43.487.32"
$ws.Range("E18").Value = "  +4.64%  "

$ws.Range("D19").Value = "0.0₃0931"
$ws.Range("E19").Value = "  +4.42%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.33"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").Value = "  +3.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "249.81"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("E24").Value = "  +7.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +2.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.86"
$ws.Range("E26").Value = "  +3.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.41"
$ws.Range("E27").Value = "  +1.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.142"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.60"
$ws.Range("E29").Value = "  +3.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("E30").Value = "  +5.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.04"
$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  +2.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0659"
$ws.Range("E35").Value = "  +5.46%  "

$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.56"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.44"
$ws.Range("E37").Value = "  +3.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("E39").Value = "  +4.71%  "

$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.84"
$ws.Range("E41").Value = "  +2.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.65"
$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("E43").Value = "  -12.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0972"
$ws.Range("E44").Value = "  -0.57%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.37"
$ws.Range("E46").Value = "  -0.29%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.474.85"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.82"
$ws.Range("E48").Value = "  +1.96%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  +10.50%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.09"
$ws.Range("E50").Value = "  +1.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.78"
$ws.Range("E51").Value = "  -1.09%  "
